$d = $word.ActiveDocument
$edits = 0

# --- Edit for paragraph containing: Siemes s4028198 ---
$rng = $d.Content
$found = $rng.Find.Execute('Siemes s4028198', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text: Siemes s4028198"
}
$target = $rng.Duplicate
$target.Expand(4)
$xml0 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="19DB0969" w14:textId="020E6571" w:rsidR="00C844BC" w:rsidRDefault="00FE585B"><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r w:rsidRPr="00C75A49"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:lang w:val="nl-NL"/></w:rPr><w:t>Achtergrond artikel</w:t></w:r><w:r w:rsidR="00C75A49" w:rsidRPr="00C75A49"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> – Jort </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:lang w:val="nl-NL"/></w:rPr><w:t>Siemes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> s4028198</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml0)
$edits = $edits + 1

# --- Edit for paragraph containing: oorlog uit op ons continent ---
$rng = $d.Content
$found = $rng.Find.Execute('oorlog uit op ons continent', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text: oorlog uit op ons continent"
}
$target = $rng.Duplicate
$target.Expand(4)
$xml10 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1D6D308D" w14:textId="3222EDDC" w:rsidR="003337A2" w:rsidRPr="003337A2" w:rsidRDefault="003337A2" w:rsidP="003337A2"><w:pPr><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r w:rsidRPr="003337A2"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>In 2019 werd de Europese Unie na een langdurige periode van "gratis geld" plotseling geconfronteerd met een enorme schok: de coronapandemie. Hoewel de schokgolven door de wereldeconomie hard werden gevoeld, reageerden centrale banken in eerste instantie terughoudend. Pas toen ze begonnen in te grijpen met steunpakketten voor de zwaarst getroffen sectoren, begonnen ze hun beleid van massa</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">al geld printen </w:t></w:r><w:r w:rsidRPr="003337A2"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>toe te passen om de immense kosten van deze steunpakketten te dekken. Dit, in combinatie met tekorten bij bijvoorbeeld computerchipfabrikanten, resulteerde in snel stijgende inflatie. Als klap op de vuurpijl brak er ook nog een oorlog</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> tussen Rusland en Oekraïne</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> uit op ons continent, wat leidde tot een ongekende stijging van de energieprijzen.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml10)
$edits = $edits + 1

# --- Edit for paragraph containing: Mujagic ---
$rng = $d.Content
$found = $rng.Find.Execute('Mujagic', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text: Mujagic"
}
$target = $rng.Duplicate
$target.Expand(4)
$xml11 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="06C04E6A" w14:textId="3FB6978C" w:rsidR="003337A2" w:rsidRPr="003337A2" w:rsidRDefault="003337A2" w:rsidP="003337A2"><w:pPr><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r w:rsidRPr="003337A2"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Economen die </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>eerst</w:t></w:r><w:r w:rsidRPr="003337A2"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> niet kritisch waren over de dagelijkse geldhoeveelheden die werden bijgedrukt, begonnen nu vragen te stellen over het terughoudende beleid van de Europese Centrale Bank met betrekking tot het verhogen van de </w:t></w:r><w:r w:rsidRPr="003337A2"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">rentetarieven. De meerderheid van de economen, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>eens</w:t></w:r><w:r w:rsidRPr="003337A2"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> met </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">*bron Edin </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>Mujagic</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>*</w:t></w:r><w:r w:rsidRPr="003337A2"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>, is het erover eens dat de ECB te laat is begonnen met dit beleid.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml11)
$edits = $edits + 1

# --- Edit for paragraph containing: retail expert en consultant ---
$rng = $d.Content
$found = $rng.Find.Execute('retail expert en consultant', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text: retail expert en consultant"
}
$target = $rng.Duplicate
$target.Expand(4)
$xml16 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="209DFEFF" w14:textId="3AE5B0D0" w:rsidR="00690C21" w:rsidRPr="00690C21" w:rsidRDefault="00690C21" w:rsidP="00690C21"><w:pPr><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r w:rsidRPr="00690C21"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">De inflatie, met name in de supermarkten, is al bekend bij de meeste Nederlanders. Volgens </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>retail</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> expert en consultant </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">*bron </w:t></w:r><w:r w:rsidRPr="00690C21"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Marco </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>Kesteloo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>*</w:t></w:r><w:r w:rsidRPr="00690C21"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">, is er echter vanuit de winkel zelf ook veel gaande. Een tekort aan personeel, gecombineerd met hoge energiekosten en </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">harde </w:t></w:r><w:r w:rsidRPr="00690C21"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">concurrentie, heeft geleid tot een </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>flinke</w:t></w:r><w:r w:rsidRPr="00690C21"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> stijging van de prijzen. Dit snelle tempo heeft ertoe geleid dat consumenten sommige producten helemaal niet meer kopen. Dit vormt een risico van inflatie voor producenten, aangezien klanten overstappen naar goedkopere concurrenten of essentiële producten overslaan.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml16)
$edits = $edits + 1

# --- Edit for paragraph containing: krimpflatie ---
$rng = $d.Content
$found = $rng.Find.Execute('krimpflatie', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text: krimpflatie"
}
$target = $rng.Duplicate
$target.Expand(4)
$xml17 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1CC45D2F" w14:textId="65659762" w:rsidR="00690C21" w:rsidRPr="00690C21" w:rsidRDefault="00690C21" w:rsidP="00690C21"><w:pPr><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r w:rsidRPr="00690C21"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Een oplossing die door oplettende consumenten is opgemerkt, is het fenomeen </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00690C21"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>krimpflatie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>”</w:t></w:r><w:r w:rsidRPr="00690C21"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">. De prijzen zijn dermate gestegen dat het economisch gezien niet meer zinvol is om ze verder te verhogen. De prijs blijft nu hetzelfde, maar producenten verminderen de inhoud. Dit gebeurt bijna altijd niet transparant, waardoor mensen </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">niet </w:t></w:r><w:r w:rsidRPr="00690C21"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">beseffen dat ze minder krijgen voor dezelfde prijs. Initiatieven zoals </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">*bron </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00690C21"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>foodwatch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t>*</w:t></w:r><w:r w:rsidRPr="00690C21"><w:rPr><w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> zijn in het leven geroepen om deze producten te identificeren en zowel prijsstijgingen als verminderingen in inhoud aan te tonen.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml17)
$edits = $edits + 1

Write-Output "Applied $edits edits"
